# UT clean up to xlsx files
# - Remove the "Texas Notes" sheet entirely.
# - Update the "About" sheet source/notes to reference the NREL Annual
#   Technology Baseline (2019) instead of the Interagency Working Group on
#   the Social Cost of Carbon, and drop the old whitehouse.gov hyperlink.
# - Update the "DR" sheet header label.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------
# 1. Delete the "Texas Notes" sheet (content moved/superseded elsewhere).
# ---------------------------------------------------------------------
$notes = $wb.Worksheets.Item("Texas Notes")
$notes.Delete() | Out-Null

# ---------------------------------------------------------------------
# 2. Rewrite the "About" sheet.
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Drop the old hyperlink (http://www.whitehouse.gov/omb/circulars_a094#8)
# before clearing, so no stale hyperlink entries are left behind.
foreach ($h in $about.Hyperlinks) { $h.Delete() }
$about.Cells.Clear()

$about.Range("A1").Value = "DR Discount Rate"
$about.Range("A1").Font.Bold = $true

$about.Range("A3").Value = "Source:"
$about.Range("A3").Font.Bold = $true
$about.Range("B3").Value = "NREL Annual Technology Baseline"

$about.Range("B4").Value = 2019
$about.Range("B4").HorizontalAlignment = -4131  ## xlLeft

$about.Range("B5").Value = "https://atb.nrel.gov/electricity/2019/data.html"
$about.Range("B5").Style = "Hyperlink"

$about.Range("B6").Value = 'See "WACC Calc" Tab'

$about.Range("A8").Value = "Notes:"
$about.Range("A8").Font.Bold = $true

$about.Range("A9").Value = "We use a 5.87% discount rate based on:"
$about.Range("A10").Value = "1) Feedback from Vibrant Clean Energy, who performed capacity expansion modeling exercises related to the Texas EPS project, and"
$about.Range("A11").Value = "2) the Interest rates used by the NREL Annual Technology Baseline"

$about.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Rewrite the "DR" sheet header label (rate value itself is unchanged).
# ---------------------------------------------------------------------
$dr = $wb.Worksheets.Item("DR")
$dr.Range("B1").Value = "Annual Perc (dimensionless)"
$dr.Range("B1").HorizontalAlignment = -4152  ## xlRight
$dr.Range("B1").WrapText = $true

$wb.Save() | Out-Null
